$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.370286
$ws.Range("H2").Value = 7.110858
$ws.Range("I2").Value = 0.3026841782318013
$ws.Range("J2").Value = 0.3026841782318014
$ws.Range("M2").Value = 6.322177333333333
$ws.Range("N2").Value = 18.966532
$ws.Range("O2").Value = 0.08271011762055308
$ws.Range("P2").Value = 0.08271011762055309
$ws.Range("Q2").Value = 14.98536842271733
$ws.Range("R2").Value = 134.868315804456
$ws.Range("S2").Value = 0.02503504398343274
$ws.Range("T2").Value = 0.02503504398343275
$ws.Range("G3").Value = 2.370286
$ws.Range("H3").Value = 7.110858
$ws.Range("I3").Value = 0.3026841782318013
$ws.Range("J3").Value = 0.3026841782318014
$ws.Range("O3").Value = 0.5401386314560596
$ws.Range("P3").Value = 0.5401386314560597
$ws.Range("Q3").Value = 97.86198623057
$ws.Range("R3").Value = 880.7578760751301
$ws.Range("S3").Value = 0.1634914177935272
$ws.Range("T3").Value = 0.1634914177935272
$ws.Range("G4").Value = 2.370286
$ws.Range("H4").Value = 7.110858
$ws.Range("I4").Value = 0.3026841782318013
$ws.Range("J4").Value = 0.3026841782318014
$ws.Range("M4").Value = 27.73243066666667
$ws.Range("N4").Value = 83.197292
$ws.Range("O4").Value = 0.3628105447549136
$ws.Range("P4").Value = 0.3628105447549136
$ws.Range("Q4").Value = 65.73379215517068
$ws.Range("R4").Value = 591.6041293965361
$ws.Range("S4").Value = 0.1098170115929732
$ws.Range("T4").Value = 0.1098170115929732
$ws.Range("G5").Value = 2.370286
$ws.Range("H5").Value = 7.110858
$ws.Range("I5").Value = 0.3026841782318013
$ws.Range("J5").Value = 0.3026841782318014
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.096171666666667
$ws.Range("N5").Value = 3.288515
$ws.Range("O5").Value = 0.01434070616847367
$ws.Range("P5").Value = 0.01434070616847367
$ws.Range("Q5").Value = 2.598240355096667
$ws.Range("R5").Value = 23.38416319587
$ws.Range("S5").Value = 0.004340704861868176
$ws.Range("T5").Value = 0.004340704861868177
$ws.Range("I6").Value = 0.2022126055089961
$ws.Range("J6").Value = 0.2022126055089961
$ws.Range("M6").Value = 6.322177333333333
$ws.Range("N6").Value = 18.966532
$ws.Range("O6").Value = 0.08271011762055308
$ws.Range("P6").Value = 0.08271011762055309
$ws.Range("Q6").Value = 10.01119520343511
$ws.Range("R6").Value = 90.100756830916
$ws.Range("S6").Value = 0.01672502838600757
$ws.Range("T6").Value = 0.01672502838600757
$ws.Range("I7").Value = 0.2022126055089961
$ws.Range("J7").Value = 0.2022126055089961
$ws.Range("O7").Value = 0.5401386314560596
$ws.Range("P7").Value = 0.5401386314560597
$ws.Range("S7").Value = 0.1092228400027932
$ws.Range("T7").Value = 0.1092228400027933
$ws.Range("I8").Value = 0.2022126055089961
$ws.Range("J8").Value = 0.2022126055089961
$ws.Range("M8").Value = 27.73243066666667
$ws.Range("N8").Value = 83.197292
$ws.Range("O8").Value = 0.3628105447549136
$ws.Range("P8").Value = 0.3628105447549136
$ws.Range("Q8").Value = 43.91442413453289
$ws.Range("R8").Value = 395.229817210796
$ws.Range("S8").Value = 0.07336486556102932
$ws.Range("T8").Value = 0.07336486556102934
$ws.Range("I9").Value = 0.2022126055089961
$ws.Range("J9").Value = 0.2022126055089961
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.096171666666667
$ws.Range("N9").Value = 3.288515
$ws.Range("O9").Value = 0.01434070616847367
$ws.Range("P9").Value = 0.01434070616847367
$ws.Range("Q9").Value = 1.735792584243889
$ws.Range("R9").Value = 15.622133258195
$ws.Range("S9").Value = 0.002899871559165993
$ws.Range("T9").Value = 0.002899871559165993
$ws.Range("G10").Value = 2.286703333333333
$ws.Range("H10").Value = 6.860109999999999
$ws.Range("I10").Value = 0.2920107190904054
$ws.Range("J10").Value = 0.2920107190904054
$ws.Range("M10").Value = 6.322177333333333
$ws.Range("N10").Value = 18.966532
$ws.Range("O10").Value = 0.08271011762055308
$ws.Range("P10").Value = 0.08271011762055309
$ws.Range("Q10").Value = 14.45694398205778
$ws.Range("R10").Value = 130.11249583852
$ws.Range("S10").Value = 0.02415224092242972
$ws.Range("T10").Value = 0.02415224092242972
$ws.Range("G11").Value = 2.286703333333333
$ws.Range("H11").Value = 6.860109999999999
$ws.Range("I11").Value = 0.2920107190904054
$ws.Range("J11").Value = 0.2920107190904054
$ws.Range("O11").Value = 0.5401386314560596
$ws.Range("P11").Value = 0.5401386314560597
$ws.Range("Q11").Value = 94.41110908981665
$ws.Range("R11").Value = 849.6999818083499
$ws.Range("S11").Value = 0.1577262701799914
$ws.Range("T11").Value = 0.1577262701799915
$ws.Range("G12").Value = 2.286703333333333
$ws.Range("H12").Value = 6.860109999999999
$ws.Range("I12").Value = 0.2920107190904054
$ws.Range("J12").Value = 0.2920107190904054
$ws.Range("M12").Value = 27.73243066666667
$ws.Range("N12").Value = 83.197292
$ws.Range("O12").Value = 0.3628105447549136
$ws.Range("P12").Value = 0.3628105447549136
$ws.Range("Q12").Value = 63.41584164690222
$ws.Range("R12").Value = 570.7425748221199
$ws.Range("S12").Value = 0.105944568067464
$ws.Range("T12").Value = 0.105944568067464
$ws.Range("G13").Value = 2.286703333333333
$ws.Range("H13").Value = 6.860109999999999
$ws.Range("I13").Value = 0.2920107190904054
$ws.Range("J13").Value = 0.2920107190904054
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.096171666666667
$ws.Range("N13").Value = 3.288515
$ws.Range("O13").Value = 0.01434070616847367
$ws.Range("P13").Value = 0.01434070616847367
$ws.Range("Q13").Value = 2.506619404072222
$ws.Range("R13").Value = 22.55957463665
$ws.Range("S13").Value = 0.004187639920520208
$ws.Range("T13").Value = 0.004187639920520208
$ws.Range("G14").Value = 1.590394666666667
$ws.Range("H14").Value = 4.771184
$ws.Range("I14").Value = 0.2030924971687972
$ws.Range("J14").Value = 0.2030924971687972
$ws.Range("M14").Value = 6.322177333333333
$ws.Range("N14").Value = 18.966532
$ws.Range("O14").Value = 0.08271011762055308
$ws.Range("P14").Value = 0.08271011762055309
$ws.Range("Q14").Value = 10.05475711265422
$ws.Range("R14").Value = 90.49281401388801
$ws.Range("S14").Value = 0.01679780432868306
$ws.Range("T14").Value = 0.01679780432868306
$ws.Range("G15").Value = 1.590394666666667
$ws.Range("H15").Value = 4.771184
$ws.Range("I15").Value = 0.2030924971687972
$ws.Range("J15").Value = 0.2030924971687972
$ws.Range("O15").Value = 0.5401386314560596
$ws.Range("P15").Value = 0.5401386314560597
$ws.Range("Q15").Value = 65.66261665069332
$ws.Range("R15").Value = 590.96354985624
$ws.Range("S15").Value = 0.1096981034797478
$ws.Range("T15").Value = 0.1096981034797478
$ws.Range("G16").Value = 1.590394666666667
$ws.Range("H16").Value = 4.771184
$ws.Range("I16").Value = 0.2030924971687972
$ws.Range("J16").Value = 0.2030924971687972
$ws.Range("M16").Value = 27.73243066666667
$ws.Range("N16").Value = 83.197292
$ws.Range("O16").Value = 0.3628105447549136
$ws.Range("P16").Value = 0.3628105447549136
$ws.Range("Q16").Value = 44.10550982596978
$ws.Range("R16").Value = 396.949588433728
$ws.Range("S16").Value = 0.07368409953344705
$ws.Range("T16").Value = 0.07368409953344705
$ws.Range("G17").Value = 1.590394666666667
$ws.Range("H17").Value = 4.771184
$ws.Range("I17").Value = 0.2030924971687972
$ws.Range("J17").Value = 0.2030924971687972
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 1.096171666666667
$ws.Range("N17").Value = 3.288515
$ws.Range("O17").Value = 0.01434070616847367
$ws.Range("P17").Value = 0.01434070616847367
$ws.Range("Q17").Value = 1.743345572417778
$ws.Range("R17").Value = 15.69011015176
$ws.Range("S17").Value = 0.00291248982691929
$ws.Range("T17").Value = 0.00291248982691929
